$d = $word.ActiveDocument
$n = $d.Paragraphs.Count
$lastReal = $d.Paragraphs.Item($n - 1)
$lastReal.Range.InsertParagraphAfter() | Out-Null
$n2 = $d.Paragraphs.Count
$placeholder = $d.Paragraphs.Item($n2 - 1)
$insertRange = $placeholder.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:pageBreakBefore/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Artificial Intelligence (AI) Generative Pre-Trained (GPT) Large Language Model (LLM) Prompts</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Negative</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">provide just one negative, contextual response to the question "Will I win the lottery?", speaking as a friend in </w:t>
      </w:r>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">empathetic </w:t>
      </w:r>
      <w:r>
        <w:t>|</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> apologetic</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> | </w:t>
      </w:r>
      <w:r>
        <w:t>disappointed</w:t>
      </w:r>
      <w:r>
        <w:t>}</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tone</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Neutral</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">provide just one neutral, contextual response to the question "Will I win the lottery?", speaking as a friend in </w:t>
      </w:r>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:r>
        <w:t>wishy-washy</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> | </w:t>
      </w:r>
      <w:r>
        <w:t>apathetic</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> | </w:t>
      </w:r>
      <w:r>
        <w:t>sanguine</w:t>
      </w:r>
      <w:r>
        <w:t>}</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tone</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Positive</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">provide just one </w:t>
      </w:r>
      <w:r>
        <w:t>positive</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, committed, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">contextual </w:t>
      </w:r>
      <w:r>
        <w:t>response to the question "Will I win the lottery</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">?", speaking as a friend </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">in </w:t>
      </w:r>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:r>
        <w:t>happy</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> |</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ecstatic</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> |</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> optimistic</w:t>
      </w:r>
      <w:r>
        <w:t>}</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tone</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Sample questions</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>Will I win the lottery?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>Will I pass my test?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>What is a dog? -&gt; need to validate input as Yes/No/Maybe-type question</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Validation</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>briefly respond true or false if the following could be a Magic-8-Ball type of question: "Will I win the lottery?"</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Google AI API Key</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>AIzaSyBNvpP3kciiZd0bmSoTT8zm-x4wa5Z1c54</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>curl -H "Content-Type: application/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>json</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">" -d </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>"{ ''</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">prompt'': { ''text'': </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>''</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">provide just one positive, committed, contextual response to the question </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>"Will I win the lottery?</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>", speaking as a friend in happy tone</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>''</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>} }" "https://generativelanguage.googleapis.com/v1beta2/models/text-bison-001:generateText?key=AIzaSyBNvpP3kciiZd0bmSoTT8zm-x4wa5Z1c54"</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xml) | Out-Null
